$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F3").Value = 1026.54
$wsSummary.Activate()
$wsSummary.Range("B4").Select()

# --- Repayment Schedule sheet (selection only) ---
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsRepay.Activate()
$wsRepay.Range("B6").Select()

# --- Transactions sheet ---
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Range("A2").Value = 214
$wsTrans.Range("A3").Value = 212
$wsTrans.Activate()
$wsTrans.Range("C3").Select()
